$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) cells to remain text so numeric-looking values
# (e.g. "20.00", "0.999", "2.20", "0.0000171") keep their exact literal
# formatting instead of being auto-converted to floating point numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated values
$ws.Range("D2").Value = "60.956.24"
$ws.Range("E2").Value = "  +0.98%  "
$ws.Range("D3").Value = "2.358.32"
$ws.Range("E3").Value = "  +1.08%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "547.11"
$ws.Range("E5").Value = "  +1.20%  "
$ws.Range("D6").Value = "138.75"
$ws.Range("E6").Value = "  +1.82%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  -0.67%  "
$ws.Range("D9").Value = "2.357.82"
$ws.Range("E9").Value = "  +1.06%  "
$ws.Range("E10").Value = "  +2.76%  "
$ws.Range("E11").Value = "  +1.72%  "
$ws.Range("D12").Value = "5.32"
$ws.Range("E12").Value = "  +1.04%  "
$ws.Range("E13").Value = "  +3.03%  "
$ws.Range("D14").Value = "25.42"
$ws.Range("E14").Value = "  +4.15%  "
$ws.Range("D15").Value = "0.0000171"
$ws.Range("E15").Value = "  +7.18%  "
$ws.Range("D16").Value = "2.784.14"
$ws.Range("E16").Value = "  +1.09%  "
$ws.Range("D17").Value = "61.173.99"
$ws.Range("E17").Value = "  +1.59%  "
$ws.Range("D18").Value = "2.362.10"
$ws.Range("E18").Value = "  +1.27%  "
$ws.Range("D19").Value = "10.93"
$ws.Range("E19").Value = "  +4.36%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "320.22"
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D21").Value = "4.13"
$ws.Range("E21").Value = "  +2.14%  "
$ws.Range("D22").Value = "6.59"
$ws.Range("E22").Value = "  +0.76%  "
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").Value = "64.06"
$ws.Range("E24").Value = "  +1.96%  "
$ws.Range("D25").Value = "1.71"
$ws.Range("E25").Value = "  -7.45%  "
$ws.Range("D26").Value = "8.87"
$ws.Range("E26").Value = "  +3.44%  "
$ws.Range("E27").Value = "  +0.63%  "
$ws.Range("D28").Value = "2.478.05"
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("D29").Value = "526.81"
$ws.Range("E29").Value = "  +5.68%  "
$ws.Range("D30").Value = "8.16"
$ws.Range("E30").Value = "  +3.08%  "
$ws.Range("D31").Value = "0.0₃0899"
$ws.Range("E31").Value = "  +1.28%  "
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("D33").Value = "0.146"
$ws.Range("E33").Value = "  +1.42%  "
$ws.Range("D34").Value = "1.83"
$ws.Range("E34").Value = "  +2.84%  "
$ws.Range("E35").Value = "  -1.24%  "
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("E37").Value = "  +7.22%  "
$ws.Range("E38").Value = "  +1.77%  "
$ws.Range("E39").Value = "  +5.11%  "
$ws.Range("E40").Value = "  +1.89%  "
$ws.Range("D41").Value = "18.41"
$ws.Range("E41").Value = "  +1.03%  "
$ws.Range("D42").Value = "145.30"
$ws.Range("E42").Value = "  +5.66%  "
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").Value = "41.43"
$ws.Range("E44").Value = "  +3.47%  "
$ws.Range("D45").Value = "146.45"
$ws.Range("E45").Value = "  +4.04%  "
$ws.Range("D46").Value = "2.20"
$ws.Range("E46").Value = "  +4.96%  "
$ws.Range("D47").Value = "3.58"
$ws.Range("E47").Value = "  +1.75%  "
$ws.Range("D48").Value = "0.0525"
$ws.Range("E48").Value = "  +3.48%  "
$ws.Range("D49").Value = "20.00"
$ws.Range("E49").Value = "  +3.60%  "
$ws.Range("D50").Value = "0.577"
$ws.Range("E50").Value = "  +2.07%  "
$ws.Range("D51").Value = "0.0899"
$ws.Range("E51").Value = "  +0.53%  "
